# more options for tweaking sigmas etc
#
# Adds 7 new output columns (Diad2_refit, HB1_Cent, HB1_Area, HB2_Cent,
# HB2_Area, C13_Cent, C13_Area) to the fitting results sheet, replaces the
# old boolean "refit" flag columns (M) with a textual warnings flag, and
# refreshes the numeric fit results with values from the re-run fitting
# routine (tweaked sigma handling).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new header cells (W1:AC1), matching style of existing header row ---
$ws.Range("V1").Copy()
$ws.Range("W1:AC1").PasteSpecial(-4122)
$ws.Range("W1").Value = "Diad2_refit"
$ws.Range("X1").Value = "HB1_Cent"
$ws.Range("Y1").Value = "HB1_Area"
$ws.Range("Z1").Value = "HB2_Cent"
$ws.Range("AA1").Value = "HB2_Area"
$ws.Range("AB1").Value = "C13_Cent"
$ws.Range("AC1").Value = "C13_Area"

# --- Row 2 ---
$ws.Range("C2").Value = 104.4416826061054
$ws.Range("D2").Value = 1282.731560903755
$ws.Range("E2").Value = 14938.35919160645
$ws.Range("F2").Value = 1282.731610906255
$ws.Range("G2").Value = 46056.9278321534
$ws.Range("H2").Value = 1.10685943665716
$ws.Range("J2").Value = 26.05434347263179
$ws.Range("K2").Value = 0.7347758307474805
$ws.Range("L2").Value = 2.21371887331432
$ws.Range("M2").Value = "Flagged Warnings:"
$ws.Range("W2").Value = "Flagged Warnings:"
$ws.Range("X2").Value = 1262.700591628001
$ws.Range("Y2").Value = 6502.745849278953
$ws.Range("Z2").Value = 1408.426449771047
$ws.Range("AA2").Value = 9658.957575289749
$ws.Range("AB2").Value = 1369.639259871846
$ws.Range("AC2").Value = 1270.251510690677

# --- Row 3 ---
$ws.Range("C3").Value = 104.9482015439476
$ws.Range("D3").Value = 1281.782458763857
$ws.Range("E3").Value = 19888.96018268805
$ws.Range("F3").Value = 1281.782408761357
$ws.Range("G3").Value = 59002.26018971133
$ws.Range("H3").Value = 1.071727856268716
$ws.Range("J3").Value = 27.54772052758774
$ws.Range("K3").Value = 0.7429253089663556
$ws.Range("L3").Value = 2.143455712537432
$ws.Range("M3").Value = "Flagged Warnings:"
$ws.Range("W3").Value = "Flagged Warnings:"
$ws.Range("X3").Value = 1261.760446403409
$ws.Range("Y3").Value = 8678.34827550104
$ws.Range("Z3").Value = 1408.059745170119
$ws.Range("AA3").Value = 10528.49189397234
$ws.Range("AB3").Value = 1369.532303367029
$ws.Range("AC3").Value = 1329.214833356688

# --- Row 4 ---
$ws.Range("C4").Value = 104.5284597401205
$ws.Range("D4").Value = 1282.556569148814
$ws.Range("E4").Value = 26070.74957900919
$ws.Range("F4").Value = 1282.556419141314
$ws.Range("G4").Value = 77582.4028390746
$ws.Range("H4").Value = 1.096548351424453
$ws.Range("J4").Value = 36.40482454093567
$ws.Range("K4").Value = 0.6874518354469279
$ws.Range("L4").Value = 2.193096702848907
$ws.Range("M4").Value = "Flagged Warnings:"
$ws.Range("W4").Value = "Flagged Warnings:"
$ws.Range("X4").Value = 1262.500374127099
$ws.Range("Y4").Value = 11118.86746785118
$ws.Range("Z4").Value = 1408.440059140436
$ws.Range("AA4").Value = 14245.34123706916
$ws.Range("AB4").Value = 1369.606721625361
$ws.Range("AC4").Value = 1726.645932763279

# --- Row 5 ---
$ws.Range("C5").Value = 104.5928477781515
$ws.Range("D5").Value = 1282.422380045482
$ws.Range("E5").Value = 27783.21832541282
$ws.Range("F5").Value = 1282.422230037982
$ws.Range("G5").Value = 82189.29412392301
$ws.Range("H5").Value = 1.086462323972746
$ws.Range("J5").Value = 33.23821993934605
$ws.Range("K5").Value = 0.6975720889834258
$ws.Range("L5").Value = 2.172924647945491
$ws.Range("M5").Value = "Flagged Warnings:"
$ws.Range("W5").Value = "Flagged Warnings: G_HighAmp"
$ws.Range("X5").Value = 1262.422595420142
$ws.Range("Y5").Value = 12082.02179617755
$ws.Range("Z5").Value = 1408.403012072601
$ws.Range("AA5").Value = 7959.804678916306
$ws.Range("AB5").Value = 1370.524840691275
$ws.Range("AC5").Value = 631.2480053987308

# --- Row 6 ---
$ws.Range("C6").Value = 104.6494847777587
$ws.Range("D6").Value = 1282.314775778605
$ws.Range("E6").Value = 28640.20313606597
$ws.Range("F6").Value = 1282.314625771104
$ws.Range("G6").Value = 84494.64331070393
$ws.Range("H6").Value = 1.082228038054093
$ws.Range("J6").Value = 35.00185500506925
$ws.Range("K6").Value = 0.7019658137253
$ws.Range("L6").Value = 2.164456076108186
$ws.Range("M6").Value = "Flagged Warnings:"
$ws.Range("W6").Value = "Flagged Warnings:"
$ws.Range("X6").Value = 1262.267393849275
$ws.Range("Y6").Value = 12397.33240083351
$ws.Range("Z6").Value = 1408.307848334188
$ws.Range("AA6").Value = 15510.18480839576
$ws.Range("AB6").Value = 1369.588992614416
$ws.Range("AC6").Value = 1909.746848901286

# --- Row 7 ---
$ws.Range("C7").Value = 104.7006071670767
$ws.Range("D7").Value = 1282.221380844187
$ws.Range("E7").Value = 28849.94416119181
$ws.Range("F7").Value = 1282.221330841687
$ws.Range("G7").Value = 84972.64515497365
$ws.Range("H7").Value = 1.075852168341284
$ws.Range("J7").Value = 35.13331555840349
$ws.Range("K7").Value = 0.7121180399087974
$ws.Range("L7").Value = 2.151704336682569
$ws.Range("M7").Value = "Flagged Warnings:"
$ws.Range("W7").Value = "Flagged Warnings:"
$ws.Range("X7").Value = 1262.172418639208
$ws.Range("Y7").Value = 12246.49910889012
$ws.Range("Z7").Value = 1408.259277165743
$ws.Range("AA7").Value = 16029.08321640206
$ws.Range("AB7").Value = 1369.568946009335
$ws.Range("AC7").Value = 1853.023849203802

# --- Row 8 ---
$ws.Range("C8").Value = 104.8871737837685
$ws.Range("D8").Value = 1281.888681855408
$ws.Range("E8").Value = 30295.70622960277
$ws.Range("F8").Value = 1281.888531847908
$ws.Range("G8").Value = 89463.14199920882
$ws.Range("H8").Value = 1.071075902070278
$ws.Range("J8").Value = 39.4268050741318
$ws.Range("K8").Value = 0.733267951110507
$ws.Range("L8").Value = 2.142151804140556
$ws.Range("M8").Value = "Flagged Warnings:"
$ws.Range("W8").Value = "Flagged Warnings:"
$ws.Range("X8").Value = 1261.867346691758
$ws.Range("Y8").Value = 13265.06820181278
$ws.Range("Z8").Value = 1408.092802551151
$ws.Range("AA8").Value = 16137.34347282708
$ws.Range("AB8").Value = 1369.543265961725
$ws.Range("AC8").Value = 2004.488423201079

# --- Row 9 ---
$ws.Range("C9").Value = 104.8021260425749
$ws.Range("D9").Value = 1282.036649610016
$ws.Range("E9").Value = 29796.23314190377
$ws.Range("F9").Value = 1282.036499602515
$ws.Range("G9").Value = 87724.00063667075
$ws.Range("H9").Value = 1.073059888351807
$ws.Range("J9").Value = 38.25837991113856
$ws.Range("K9").Value = 0.7205048029720834
$ws.Range("L9").Value = 2.146119776703614
$ws.Range("M9").Value = "Flagged Warnings:"
$ws.Range("W9").Value = "Flagged Warnings:"
$ws.Range("X9").Value = 1262.019647960024
$ws.Range("Y9").Value = 12918.63043263268
$ws.Range("Z9").Value = 1408.182013521922
$ws.Range("AA9").Value = 15827.56809427078
$ws.Range("AB9").Value = 1369.550184344768
$ws.Range("AC9").Value = 1955.040423331327
